$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Existing hyperlinks (F2:F5) get rebuilt below together with the three
# brand-new ones (F6:F8), so clear the collection first and re-add in
# order -> rId1..rId7 land on the same cells the diff expects.
$ws.Hyperlinks.Delete()

# Row 2
$ws.Cells.Item(2, 1).Value = '2026-01-05 01:44:10'
$ws.Cells.Item(2, 2).Value = '【急募】医療診断AIの深層学習モデル開発の専門家募集'
$ws.Cells.Item(2, 3).Value = 'システム開発'
$ws.Cells.Item(2, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = '期限情報なし'
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5464587'
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), 'https://www.lancers.jp/work/detail/5464587')
$ws.Cells.Item(2, 6).Style = "Hyperlink"
$ws.Cells.Item(2, 7).Value = 375
$ws.Cells.Item(2, 8).Value = '🔥AI,Ai ◆開発'

# Row 3
$ws.Cells.Item(3, 1).Value = '2026-01-05 01:44:10'
$ws.Cells.Item(3, 2).Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Cells.Item(3, 3).Value = 'システム開発'
$ws.Cells.Item(3, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = '期限情報なし'
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), 'https://www.lancers.jp/work/detail/5217096')
$ws.Cells.Item(3, 6).Style = "Hyperlink"
$ws.Cells.Item(3, 7).Value = 243
$ws.Cells.Item(3, 8).Value = '🔥API ◆ツール'

# Row 4
$ws.Cells.Item(4, 1).Value = '2026-01-05 01:44:10'
$ws.Cells.Item(4, 2).Value = '【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています'
$ws.Cells.Item(4, 3).Value = 'システム開発'
$ws.Cells.Item(4, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = '期限情報なし'
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5405023'
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), 'https://www.lancers.jp/work/detail/5405023')
$ws.Cells.Item(4, 6).Style = "Hyperlink"
$ws.Cells.Item(4, 7).Value = 178
$ws.Cells.Item(4, 8).Value = '★bot ◆ツール'

# Row 5
$ws.Cells.Item(5, 1).Value = '2026-01-05 01:44:10'
$ws.Cells.Item(5, 2).Value = '【急募】kintone案件管理アプリにExcel見積計算式組込'
$ws.Cells.Item(5, 3).Value = 'システム開発'
$ws.Cells.Item(5, 4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = '期限情報なし'
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5464763'
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), 'https://www.lancers.jp/work/detail/5464763')
$ws.Cells.Item(5, 6).Style = "Hyperlink"
$ws.Cells.Item(5, 7).Value = 55
$ws.Cells.Item(5, 8).Value = '◇アプリ'

# Row 6
$ws.Cells.Item(6, 1).Value = '2026-01-05 01:44:10'
$ws.Cells.Item(6, 2).Value = '【急募】メール問い合わせ時の自動SMS送信システム構築'
$ws.Cells.Item(6, 3).Value = 'システム開発'
$ws.Cells.Item(6, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = '期限情報なし'
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5464796'
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), 'https://www.lancers.jp/work/detail/5464796')
$ws.Cells.Item(6, 6).Style = "Hyperlink"
$ws.Cells.Item(6, 7).Value = 33

# Row 7
$ws.Cells.Item(7, 1).Value = '2026-01-05 01:44:10'
$ws.Cells.Item(7, 2).Value = '金融機関の入出金伝票印刷システム構築依頼'
$ws.Cells.Item(7, 3).Value = 'システム開発'
$ws.Cells.Item(7, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = '期限情報なし'
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5464833'
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), 'https://www.lancers.jp/work/detail/5464833')
$ws.Cells.Item(7, 6).Style = "Hyperlink"
$ws.Cells.Item(7, 7).Value = 28

# Row 8
$ws.Cells.Item(8, 1).Value = '2026-01-05 01:44:10'
$ws.Cells.Item(8, 2).Value = '【急募】簡単なHP作成とAWS構築をしてくれる方募集'
$ws.Cells.Item(8, 3).Value = 'システム開発'
$ws.Cells.Item(8, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = '期限情報なし'
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5457524'
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), 'https://www.lancers.jp/work/detail/5457524')
$ws.Cells.Item(8, 6).Style = "Hyperlink"
$ws.Cells.Item(8, 7).Value = 18
